# Update the "dSF" (column F) values for the rows whose source data was
# re-pulled, as part of the repull/mean-calculation refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = 3
$ws.Range("F6").Value = -2
$ws.Range("F9").Value = 1
$ws.Range("F11").Value = -2
$ws.Range("F13").Value = 0
